$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Runmode column (E4:E8) from "No" to "Yes"
$ws.Range("E4:E8").Value = "Yes"

# Update the active selection to reflect the new focus area (E2:E8)
$ws.Range("E2:E8").Select()
